$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force columns D and E to text format to avoid Excel auto-numeric/percent conversion
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '21.214.17'
$ws.Range("E2").Value = '  +3.89%  '
$ws.Range("D3").Value = '1.539.34'
$ws.Range("E3").Value = '  +5.21%  '
$ws.Range("D4").Value = '1.001'
$ws.Range("E4").Value = '  -0.64%  '
$ws.Range("D5").Value = '0.9592'
$ws.Range("E5").Value = '  +1.17%  '
$ws.Range("D6").Value = '281.44'
$ws.Range("E6").Value = '  +2.50%  '
$ws.Range("E7").Value = '  -0.84%  '
$ws.Range("E8").Value = '  +3.22%  '
$ws.Range("D9").Value = '40.80'
$ws.Range("E9").Value = '  +2.99%  '
$ws.Range("D10").Value = '1.104'
$ws.Range("E10").Value = '  +6.57%  '
$ws.Range("D11").Value = '0.06791'
$ws.Range("E11").Value = '  +3.47%  '
$ws.Range("D12").Value = '0.9955'
$ws.Range("E12").Value = '  -0.35%  '
$ws.Range("B13").Value = 'Polkadot'
$ws.Range("C13").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D13").Value = '5.651'
$ws.Range("E13").Value = '  +4.48%  '
$ws.Range("B14").Value = 'Solana'
$ws.Range("C14").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range("D14").Value = '18.74'
$ws.Range("E14").Value = '  +4.28%  '
$ws.Range("D15").Value = '6.338'
$ws.Range("E15").Value = '  +3.57%  '
$ws.Range("D16").Value = '0.00001044'
$ws.Range("E16").Value = '  +1.89%  '
$ws.Range("D17").Value = '0.9603'
$ws.Range("E17").Value = '  -0.44%  '
$ws.Range("D18").Value = '1.531.51'
$ws.Range("E18").Value = '  +4.88%  '
$ws.Range("D19").Value = '0.06035'
$ws.Range("E19").Value = '  +4.64%  '
$ws.Range("D20").Value = '71.79'
$ws.Range("E20").Value = '  +3.19%  '
$ws.Range("D21").Value = '5.662'
$ws.Range("E21").Value = '  +4.45%  '
$ws.Range("D22").Value = '15.01'
$ws.Range("D23").Value = '11.34'
$ws.Range("E23").Value = '  +4.36%  '
$ws.Range("D24").Value = '2.304'
$ws.Range("E24").Value = '  +3.05%  '
$ws.Range("D25").Value = '21.269.41'
$ws.Range("E25").Value = '  +4.02%  '
$ws.Range("D26").Value = '147.09'
$ws.Range("E26").Value = '  +4.04%  '
$ws.Range("D27").Value = '2.203'
$ws.Range("E27").Value = '  +5.87%  '
$ws.Range("D28").Value = '17.67'
$ws.Range("E28").Value = '  +3.23%  '
$ws.Range("D29").Value = '1.700.64'
$ws.Range("E29").Value = '  +5.40%  '
$ws.Range("D30").Value = '117.54'
$ws.Range("E30").Value = '  +5.04%  '
$ws.Range("D31").Value = '4.055'
$ws.Range("E31").Value = '  +5.19%  '
$ws.Range("D32").Value = '0.8487'
$ws.Range("E32").Value = '  +7.74%  '
$ws.Range("D33").Value = '5.163'
$ws.Range("E33").Value = '  +5.94%  '
$ws.Range("D34").Value = '0.08003'
$ws.Range("E34").Value = '  +2.60%  '
$ws.Range("D35").Value = '1.498'
$ws.Range("E35").Value = '  -0.72%  '
$ws.Range("D36").Value = '1.219'
$ws.Range("E36").Value = '  +7.62%  '
$ws.Range("D37").Value = '4.940'
$ws.Range("E37").Value = '  +6.01%  '
$ws.Range("D38").Value = '0.05847'
$ws.Range("E38").Value = '  +2.68%  '
$ws.Range("D39").Value = '0.02091'
$ws.Range("E39").Value = '  +3.13%  '
$ws.Range("D40").Value = '10.72'
$ws.Range("E40").Value = '  +3.80%  '
$ws.Range("D41").Value = '7.715'
$ws.Range("E41").Value = '  +4.08%  '
$ws.Range("D42").Value = '0.1910'
$ws.Range("E42").Value = '  +2.93%  '
$ws.Range("D43").Value = '0.9594'
$ws.Range("E43").Value = '  +0.68%  '
$ws.Range("D44").Value = '0.5451'
$ws.Range("E44").Value = '  +3.70%  '
$ws.Range("D45").Value = '12.47'
$ws.Range("E45").Value = '  +4.83%  '
$ws.Range("D46").Value = '3.553'
$ws.Range("E46").Value = '  +1.96%  '
$ws.Range("B47").Value = 'Decentraland'
$ws.Range("C47").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D47").Value = '0.5440'
$ws.Range("E47").Value = '  +5.89%  '
$ws.Range("B48").Value = 'Quant'
$ws.Range("C48").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D48").Value = '121.49'
$ws.Range("E48").Value = '  +3.85%  '
$ws.Range("D49").Value = '1.867'
$ws.Range("E49").Value = '  +6.87%  '
$ws.Range("D50").Value = '0.06622'
$ws.Range("E50").Value = '  +3.15%  '
$ws.Range("D51").Value = '69.94'
$ws.Range("E51").Value = '  +5.80%  '

# Restore default style (remove the temporary text-number-format xf) to match original formatting
$ws.Range("D2:E51").Style = "Normal"
